$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 4200.2
$ws.Range("J51").Value = 4750
$ws.Range("L51").Value = 4750
$ws.Range("N51").Value = -5718

# Row 100
$ws.Range("H100").Value = 3313.7368
$ws.Range("I100").Value = 1944.8462
$ws.Range("K100").Value = 1944.8462
$ws.Range("M100").Value = -1403.8462

# Row 128
$ws.Range("H128").Value = 81495.375
$ws.Range("J128").Value = 81495.375
$ws.Range("L128").Value = 81495.375
$ws.Range("N128").Value = -91455.375

# Row 137
$ws.Range("H137").Value = 12846.444
$ws.Range("I137").Value = 1484.6364
$ws.Range("J137").Value = 17845.64
$ws.Range("K137").Value = 4453.9092
$ws.Range("L137").Value = 53536.92
$ws.Range("M137").Value = -1903.9092
$ws.Range("N137").Value = -58636.92

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1134.4127
$ws.Range("I32").Value = 1142.8167
$ws.Range("K32").Value = 1142.8167
$ws.Range("M32").Value = -855.8167000000001

# Row 74
$ws.Range("H74").Value = 262762.3
$ws.Range("I74").Value = 278825.06
$ws.Range("K74").Value = 278825.06
$ws.Range("M74").Value = -277951.06

# Row 77
$ws.Range("H77").Value = 262762.3
$ws.Range("I77").Value = 278825.06
$ws.Range("K77").Value = 1394125.3
$ws.Range("M77").Value = -1389757.3

# Row 97
$ws.Range("H97").Value = 1810.1305
$ws.Range("J97").Value = 1072.4286
$ws.Range("L97").Value = 1072.4286
$ws.Range("N97").Value = -2064.4286

# Row 132
$ws.Range("H132").Value = 199712.75
$ws.Range("I132").Value = 243643.88
$ws.Range("J132").Value = 16000.818
$ws.Range("K132").Value = 730931.64
$ws.Range("L132").Value = 48002.454
$ws.Range("M132").Value = -728401.64
$ws.Range("N132").Value = -53062.454

# Row 139
$ws.Range("H139").Value = 59214.285
$ws.Range("J139").Value = 59214.285
$ws.Range("L139").Value = 59214.285
$ws.Range("N139").Value = -69494.285

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2351.75
$ws.Range("I20").Value = 1736.6
$ws.Range("K20").Value = 1736.6
$ws.Range("M20").Value = -1489.6

# Row 22
$ws.Range("H22").Value = 1001
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# Row 86
$ws.Range("H86").Value = 2707.111
$ws.Range("I86").Value = 2545.625
$ws.Range("K86").Value = 2545.625
$ws.Range("M86").Value = -1422.625

# Row 89
$ws.Range("H89").Value = 2707.111
$ws.Range("I89").Value = 2545.625
$ws.Range("K89").Value = 12728.125
$ws.Range("M89").Value = -7112.125

# Row 105
$ws.Range("H105").Value = 27779622
$ws.Range("I105").Value = 2211.8
$ws.Range("K105").Value = 2211.8
$ws.Range("M105").Value = -464.8000000000002

# Row 134
$ws.Range("H134").Value = 49761.832
$ws.Range("I134").Value = 2846.6667
$ws.Range("J134").Value = 77910.92999999999
$ws.Range("K134").Value = 8540.000100000001
$ws.Range("L134").Value = 233732.79
$ws.Range("M134").Value = -6005.000100000001
$ws.Range("N134").Value = -238802.79

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1218.375
$ws.Range("I16").Value = 1234
$ws.Range("J16").Value = 1202.75
$ws.Range("K16").Value = 1234
$ws.Range("L16").Value = 1202.75
$ws.Range("M16").Value = -947
$ws.Range("N16").Value = -1776.75

# Row 99
$ws.Range("H99").Value = 5368.385
$ws.Range("I99").Value = 3685
$ws.Range("K99").Value = 3685
$ws.Range("M99").Value = -2187

# Row 113
$ws.Range("H113").Value = 1218.375
$ws.Range("I113").Value = 1234
$ws.Range("J113").Value = 1202.75
$ws.Range("K113").Value = 1234
$ws.Range("L113").Value = 1202.75
$ws.Range("M113").Value = 936
$ws.Range("N113").Value = -5542.75

# Row 126
$ws.Range("H126").Value = 5368.385
$ws.Range("I126").Value = 3685
$ws.Range("K126").Value = 11055
$ws.Range("M126").Value = -8585

# Row 132
$ws.Range("H132").Value = 4293
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 182787.73
$ws.Range("J11").Value = 988.8889
$ws.Range("L11").Value = 2966.6667
$ws.Range("N11").Value = -3246.6667

# Row 122
$ws.Range("H122").Value = 8883618
$ws.Range("J122").Value = 1108.3334
$ws.Range("L122").Value = 9975.000599999999
$ws.Range("N122").Value = -14875.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7845.1816
$ws.Range("I70").Value = 7492.4287
$ws.Range("J70").Value = 8462.5
$ws.Range("K70").Value = 7492.4287
$ws.Range("L70").Value = 8462.5
$ws.Range("M70").Value = -7222.4287
$ws.Range("N70").Value = -9002.5

# Row 73
$ws.Range("H73").Value = 7845.1816
$ws.Range("I73").Value = 7492.4287
$ws.Range("J73").Value = 8462.5
$ws.Range("K73").Value = 7492.4287
$ws.Range("L73").Value = 8462.5
$ws.Range("M73").Value = -6556.4287
$ws.Range("N73").Value = -10334.5

# Row 113
$ws.Range("H113").Value = 40423280
$ws.Range("I113").Value = 1046905
$ws.Range("J113").Value = 66674200
$ws.Range("K113").Value = 1046905
$ws.Range("L113").Value = 66674200
$ws.Range("M113").Value = -1044735
$ws.Range("N113").Value = -66678540

# Row 122
$ws.Range("H122").Value = 529099.6
$ws.Range("I122").Value = 652207
$ws.Range("J122").Value = 5893.25
$ws.Range("K122").Value = 1956621
$ws.Range("L122").Value = 17679.75
$ws.Range("M122").Value = -1954171
$ws.Range("N122").Value = -22579.75

# Row 132
$ws.Range("H132").Value = 69612.85000000001
$ws.Range("I132").Value = 33246.547
$ws.Range("J132").Value = 114060.555
$ws.Range("K132").Value = 99739.641
$ws.Range("L132").Value = 342181.665
$ws.Range("M132").Value = -97209.641
$ws.Range("N132").Value = -347241.665

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 341731.47
$ws.Range("I7").Value = 592674.9
$ws.Range("J7").Value = 13574.692
$ws.Range("K7").Value = 592674.9
$ws.Range("L7").Value = 13574.692
$ws.Range("M7").Value = -592562.9
$ws.Range("N7").Value = -13798.692

# Row 16
$ws.Range("H16").Value = 1309.8928
$ws.Range("I16").Value = 1326.4615
$ws.Range("K16").Value = 1326.4615
$ws.Range("M16").Value = -1156.4615

# Row 40
$ws.Range("H40").Value = 460185.5
$ws.Range("I40").Value = 532057.2
$ws.Range("K40").Value = 532057.2
$ws.Range("M40").Value = -531921.2

# Row 46
$ws.Range("H46").Value = 2823.4707
$ws.Range("I46").Value = 2020.75
$ws.Range("J46").Value = 4750
$ws.Range("K46").Value = 2020.75
$ws.Range("L46").Value = 4750
$ws.Range("M46").Value = -1832.75
$ws.Range("N46").Value = -5126

# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 55
$ws.Range("H55").Value = 40000388
$ws.Range("I55").Value = 188.71428
$ws.Range("J55").Value = 90909736
$ws.Range("K55").Value = 188.71428
$ws.Range("L55").Value = 90909736
$ws.Range("M55").Value = -15.71428
$ws.Range("N55").Value = -90910082

# Row 99
$ws.Range("H99").Value = 25029.5
$ws.Range("I99").Value = 25029.5
$ws.Range("K99").Value = 25029.5
$ws.Range("M99").Value = -22034.5

# Row 100
$ws.Range("H100").Value = 367499.66
$ws.Range("I100").Value = 519999.5
$ws.Range("J100").Value = 62500
$ws.Range("K100").Value = 519999.5
$ws.Range("L100").Value = 62500
$ws.Range("M100").Value = -519458.5
$ws.Range("N100").Value = -63582

# Row 122
$ws.Range("H122").Value = 615857.1
$ws.Range("I122").Value = 3686.3635
$ws.Range("J122").Value = 1577839.8
$ws.Range("K122").Value = 11059.0905
$ws.Range("L122").Value = 4733519.4
$ws.Range("M122").Value = -8609.0905
$ws.Range("N122").Value = -4738419.4

# Row 126
$ws.Range("H126").Value = 341731.47
$ws.Range("I126").Value = 592674.9
$ws.Range("J126").Value = 13574.692
$ws.Range("K126").Value = 1778024.7
$ws.Range("L126").Value = 40724.076
$ws.Range("M126").Value = -1775554.7
$ws.Range("N126").Value = -45664.076

# Row 132
$ws.Range("H132").Value = 5402.853
$ws.Range("I132").Value = 3878.1428
$ws.Range("K132").Value = 11634.4284
$ws.Range("M132").Value = -9104.428400000001

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7199.8237
$ws.Range("I62").Value = 6640
$ws.Range("K62").Value = 6640
$ws.Range("M62").Value = -6016

# Row 65
$ws.Range("H65").Value = 7199.8237
$ws.Range("I65").Value = 6640
$ws.Range("K65").Value = 33200
$ws.Range("M65").Value = -30080

# Row 96
$ws.Range("H96").Value = 74413.21000000001
$ws.Range("J96").Value = 4054.3333
$ws.Range("L96").Value = 4054.3333
$ws.Range("N96").Value = -6800.3333

# Row 132
$ws.Range("H132").Value = 54140.82
$ws.Range("I132").Value = 18322.65
$ws.Range("J132").Value = 143686.25
$ws.Range("K132").Value = 54967.95
$ws.Range("L132").Value = 431058.75
$ws.Range("M132").Value = -52437.95
$ws.Range("N132").Value = -436118.75
